# Add a new spelling variant for Vietnam ("Democratic Republic Of Vietnam")
# to the Codes crosswalk table.
#
# The "Codes" sheet lists one row per (country-code, alternate spelling)
# pair; all the Vietnam (VNM) rows are identical except for the text in
# column C. This inserts a new row right after the existing "Vietnam"
# row (row 487), copies that row's values across, and then overwrites
# column C with the new spelling - pushing every following row down by
# one (506 data rows -> 507).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codes")

$sourceRow = 487
$newRow = 488
$lastCol = 25  # A..Y

$ws.Rows("$newRow`:$newRow").Insert()

for ($col = 1; $col -le $lastCol; $col++) {
    $srcVal = $ws.Cells.Item($sourceRow, $col).Value()
    if ($srcVal -ne $null -and $srcVal -ne "") {
        $ws.Cells.Item($newRow, $col).Value = $srcVal
    }
}

$ws.Cells.Item($newRow, 3).Value = "Democratic Republic Of Vietnam"

# Keep the sheet's remembered sort range in sync with the new data extent
# (A2:Y506 -> A2:Y507). The existing data is already in the stored sort
# order, so re-applying it is a no-op for row order but refreshes the
# cached range.
$lastRow = 507
$sortObj = $ws.Sort
$sortObj.SetRange($ws.Range("A1:Y" + $lastRow))
$sortObj.Header = 1
$sortObj.Apply()
